# Update column F ("dSF") values on Sheet1 to reflect repulled data / mean calc.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 3
    3  = -3
    5  = -1
    6  = -3
    7  = 2
    8  = 0
    9  = -1
    10 = 7
    11 = -1
    12 = -1
    13 = -1
    17 = -8
    18 = -2
    19 = 4
    20 = -5
    21 = 2
    24 = 3
    25 = -1
    27 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
